# Changing the forecast models for Kahraman
# Shift the timestamp column (A) by +2 days for all data rows (2..97),
# and update the production values (B) for rows 2..51 to the new forecast.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift timestamps in column A (rows 2 to 97) forward by 2 days.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 2
}

# New "Actual Production (MW)" values for rows 2 to 51.
$newB = @(833,855,869,868,899,883,859,868,842,822,845,859,856,839,816,810,816,800,785,803,807,814,816,792,747,737,717,646,535,451,396,378,349,336,337,339,354,376,390,402,366,327,283,253,222,190,172,153,142,126)

for ($i = 0; $i -lt $newB.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $newB[$i]
}
